# Automatische test-sync: 2025-08-05 18:35:50
#
# 1) Logs sheet: append a new row 35 (Testmail #14 - CE certificates).
# 2) Dashboard sheet: swap rows 3/4 category labels, append new row 9
#    (Kwaliteit / Certificaten) with count 1.
# 3) Dashboard chart: extend the category/value series references from
#    row 8 to row 9.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Logs sheet - add row 35
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(35, 1).Value = "Heb je de CE-certificaten van dit product?"
$logs.Cells.Item(35, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(35, 3).Value = "Testmail #14: Heb je de CE-certificaten van dit product?"
$logs.Cells.Item(35, 4).Value = "Kwaliteit / Certificaten"
$logs.Cells.Item(35, 5).Value = "Bedankt, we hebben dit doorgestuurd naar kwaliteit@bedrijf.nl."
$logs.Cells.Item(35, 6).Value = "2025-08-05 18:35:01"
$logs.Cells.Item(35, 7).Value = "Ja"
$logs.Cells.Item(35, 8).Value = "Ja"
$logs.Cells.Item(35, 9).Value = "Nee"
$logs.Cells.Item(35, 10).Value = "Nee"

# Extend the conditional-formatting ranges (D/G/H/I/J, rows 2-34) so they
# cover the newly-added row 35 as well.
$logs.Range("D2:D34").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D35"))
$logs.Range("G2:G34").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G35"))
$logs.Range("H2:H34").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H35"))
$logs.Range("I2:I34").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I35"))
$logs.Range("J2:J34").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J35"))

# ---------------------------------------------------------------------
# 2) Dashboard sheet - swap rows 3/4, add row 9
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Cells.Item(3, 1).Value = "Klantenservice / Contact"
$dash.Cells.Item(4, 1).Value = "Inkoop / Bestellingen"

$dash.Cells.Item(9, 1).Value = "Kwaliteit / Certificaten"
$dash.Cells.Item(9, 2).Value = 1

# ---------------------------------------------------------------------
# 3) Chart - extend series references from row 8 to row 9
# ---------------------------------------------------------------------
$chartObj = $dash.ChartObjects().Item(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES('Dashboard'!B1,'Dashboard'!`$A`$2:`$A`$9,'Dashboard'!`$B`$2:`$B`$9,1)"

Write-Output "edit complete"
